$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff / handback datetimes for row 4
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-15 08:07:04"
$wsZh.Range("G4").Value = "2016-02-15 08:07:55"

# de-de sheet: update handoff / handback datetimes for row 4
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-15 08:07:18"
$wsDe.Range("G4").Value = "2016-02-15 08:08:22"
